# Append two new rows (104, 105) of date / remn_amt data to each of the
# three worksheets, matching the upstream "Add files via upload" export.

$wb = $excel.ActiveWorkbook

# date serials for the two new rows (2025-11-06, 2025-11-07)
$date104 = 45967
$date105 = 45968

# per-sheet remn_amt values for row 104 (row 105 is always 0 across sheets)
$values104 = @(449373, 58622, 13273)

for ($i = 1; $i -le 3; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $ws.Range("A104").Value = $date104
    $ws.Range("B104").Value = $values104[$i - 1]

    $ws.Range("A105").Value = $date105
    $ws.Range("B105").Value = 0

    # Match the date/time number format used by the existing date column
    # (A2:A103) so the new cells share the same style (s="2").
    $ws.Range("A104:A105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
